$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Hoja1: update rows 2-14 (A,B,C columns; D/E/G unchanged constants) ---
$ws1.Range("A2").Value = 40317700413
$ws1.Range("B2").Value = "Sebastián Triana"
$ws1.Range("C2").Value = "Carrera 68f #96-40"
$ws1.Range("A3").Value = 40317696854
$ws1.Range("B3").Value = "Mairly Giraldo"
$ws1.Range("C3").Value = "Calle 38g sur #68c-63"
$ws1.Range("A4").Value = 40317691488
$ws1.Range("B4").Value = "Orlando Bustamante"
$ws1.Range("C4").Value = "calle 22B #44A-33"
$ws1.Range("A5").Value = 40317605804
$ws1.Range("B5").Value = "ANDRES CAMACHO"
$ws1.Range("C5").Value = "Cra 5c #48J-82sur"
$ws1.Range("A6").Value = 40317529633
$ws1.Range("B6").Value = "nicolas cruz"
$ws1.Range("C6").Value = "Transversal 13 C este #48a-09s"
$ws1.Range("A7").Value = 40317479749
$ws1.Range("B7").Value = "Sergio Hernández"
$ws1.Range("C7").Value = "Calle 57B Sur #62-13"
$ws1.Range("A8").Value = 40316841061
$ws1.Range("B8").Value = "ingrid yadira espitia castañeda"
$ws1.Range("C8").Value = "Calle 86a #69t-81"
$ws1.Range("A9").Value = 40316830437
$ws1.Range("B9").Value = "leonardo hoyos"
$ws1.Range("C9").Value = "Cra 73B #64F-58"
$ws1.Range("A10").Value = 40316455848
$ws1.Range("B10").Value = "Camilo Cuervo Díaz"
$ws1.Range("C10").Value = "Transversal 19A #95-56"
$ws1.Range("A11").Value = 40316213857
$ws1.Range("B11").Value = "STIVEN GORDILLO MURCIA"
$ws1.Range("C11").Value = "Diagonal 2A #79F-26"
$ws1.Range("A12").Value = 40316195365
$ws1.Range("B12").Value = "Camilo Romero"
$ws1.Range("C12").Value = "Carrera 3 #74A-50"
$ws1.Range("A13").Value = 40316124521
$ws1.Range("B13").Value = "Francisco Javier Mendez"
$ws1.Range("C13").Value = "Carrera90 #6a-98"
$ws1.Range("A14").Value = 100008308
$ws1.Range("B14").Value = "BRUCE CASTELLANOS"
$ws1.Range("C14").Value = "CRA 23D # 2-54 BARRIO LA FRAGÜITA"

# Hoja1: clear phone column F for rows 2 and 3 (no longer present)
$ws1.Range("F2").Value = ""
$ws1.Range("F3").Value = ""

# Hoja1: remove old rows 15-17 (data shrinks from 17 to 14 rows)
$ws1.Range("A15:H17").EntireRow.Delete()

# --- Hoja2: update rows 2-14 ---
$ws2.Range("A2").Value = 40317700413
$ws2.Range("B2").Value = 4271228271
$ws2.Range("C2").Value = "ASD00993"
$ws2.Range("D2").Value = 1
$ws2.Range("F2").Value = 44195
$ws2.Range("G2").Value = "2B"
$ws2.Range("H2").Value = "Sebastián Triana"
$ws2.Range("I2").Value = "Carrera 68f #96-40"
$ws2.Range("J2").Value = 40317700413
$ws2.Range("A3").Value = 40317696854
$ws2.Range("B3").Value = 4271221541
$ws2.Range("C3").Value = "FXXI0005"
$ws2.Range("D3").Value = 1
$ws2.Range("F3").Value = 44195
$ws2.Range("G3").Value = "6E"
$ws2.Range("H3").Value = "Mairly Giraldo"
$ws2.Range("I3").Value = "Calle 38g sur #68c-63"
$ws2.Range("J3").Value = 40317696854
$ws2.Range("A4").Value = 40317691488
$ws2.Range("B4").Value = 4271218258
$ws2.Range("C4").Value = "FXXI0010"
$ws2.Range("D4").Value = 1
$ws2.Range("F4").Value = 44195
$ws2.Range("G4").Value = "8F"
$ws2.Range("H4").Value = "Orlando Bustamante"
$ws2.Range("I4").Value = "calle 22B #44A-33"
$ws2.Range("J4").Value = 40317691488
$ws2.Range("A5").Value = 40317605804
$ws2.Range("B5").Value = 4271119539
$ws2.Range("C5").Value = "FXXI0004"
$ws2.Range("D5").Value = 1
$ws2.Range("F5").Value = 44195
$ws2.Range("G5").Value = "6E"
$ws2.Range("H5").Value = "ANDRES CAMACHO"
$ws2.Range("I5").Value = "Cra 5c #48J-82sur"
$ws2.Range("J5").Value = 40317605804
$ws2.Range("A6").Value = 40317529633
$ws2.Range("B6").Value = 4271032951
$ws2.Range("C6").Value = "ESSG0008"
$ws2.Range("D6").Value = 1
$ws2.Range("F6").Value = 44195
$ws2.Range("G6").Value = "3F"
$ws2.Range("H6").Value = "nicolas cruz"
$ws2.Range("I6").Value = "Transversal 13 C este #48a-09s"
$ws2.Range("J6").Value = 40317529633
$ws2.Range("A7").Value = 40317479749
$ws2.Range("B7").Value = 4270978937
$ws2.Range("C7").Value = "XDXI0011"
$ws2.Range("D7").Value = 1
$ws2.Range("F7").Value = 44195
$ws2.Range("G7").Value = "7F"
$ws2.Range("H7").Value = "Sergio Hernández"
$ws2.Range("I7").Value = "Calle 57B Sur #62-13"
$ws2.Range("J7").Value = 40317479749
$ws2.Range("A8").Value = 40316841061
$ws2.Range("B8").Value = 4270248701
$ws2.Range("C8").Value = "AGL01533"
$ws2.Range("D8").Value = 1
$ws2.Range("F8").Value = 44195
$ws2.Range("G8").Value = "8A"
$ws2.Range("H8").Value = "ingrid yadira espitia castañeda"
$ws2.Range("I8").Value = "Calle 86a #69t-81"
$ws2.Range("J8").Value = 40316841061
$ws2.Range("A9").Value = 40316830437
$ws2.Range("B9").Value = 4270240840
$ws2.Range("C9").Value = "076CS27184"
$ws2.Range("D9").Value = 1
$ws2.Range("F9").Value = 44195
$ws2.Range("G9").Value = "5B"
$ws2.Range("H9").Value = "leonardo hoyos"
$ws2.Range("I9").Value = "Cra 73B #64F-58"
$ws2.Range("J9").Value = 40316830437
$ws2.Range("A10").Value = 40316455848
$ws2.Range("B10").Value = 4269817874
$ws2.Range("C10").Value = "ASD00841"
$ws2.Range("D10").Value = 1
$ws2.Range("F10").Value = 44195
$ws2.Range("G10").Value = "5A"
$ws2.Range("H10").Value = "Camilo Cuervo Díaz"
$ws2.Range("I10").Value = "Transversal 19A #95-56"
$ws2.Range("J10").Value = 40316455848
$ws2.Range("A11").Value = 40316213857
$ws2.Range("B11").Value = 4269541997
$ws2.Range("C11").Value = "ST75874RS"
$ws2.Range("D11").Value = 1
$ws2.Range("F11").Value = 44195
$ws2.Range("G11").Value = "11F"
$ws2.Range("H11").Value = "STIVEN GORDILLO MURCIA"
$ws2.Range("I11").Value = "Diagonal 2A #79F-26"
$ws2.Range("J11").Value = 40316213857
$ws2.Range("A12").Value = 40316195365
$ws2.Range("B12").Value = 4269524164
$ws2.Range("C12").Value = "ACSG0001"
$ws2.Range("D12").Value = 1
$ws2.Range("F12").Value = 44195
$ws2.Range("G12").Value = "10F"
$ws2.Range("H12").Value = "Camilo Romero"
$ws2.Range("I12").Value = "Carrera 3 #74A-50"
$ws2.Range("J12").Value = 40316195365
$ws2.Range("A13").Value = 40316124521
$ws2.Range("B13").Value = 4269441491
$ws2.Range("C13").Value = "XDHW0004"
$ws2.Range("D13").Value = 1
$ws2.Range("F13").Value = 44195
$ws2.Range("G13").Value = "3G"
$ws2.Range("H13").Value = "Francisco Javier Mendez"
$ws2.Range("I13").Value = "Carrera90 #6a-98"
$ws2.Range("J13").Value = 40316124521
$ws2.Range("A14").Value = 100008308
$ws2.Range("B14").Value = 4168110107
$ws2.Range("C14").Value = "XMAP0002"
$ws2.Range("D14").Value = 1
$ws2.Range("F14").Value = 44195
$ws2.Range("G14").Value = "4D"
$ws2.Range("H14").Value = "BRUCE CASTELLANOS"
$ws2.Range("I14").Value = "CRA 23D # 2-54 BARRIO LA FRAGÜITA"
$ws2.Range("J14").Value = 100008308

# Hoja2: clear rows 15-17 (keep F-column style placeholder, matching blank template rows)
$ws2.Range("A15:E15").ClearContents()
$ws2.Range("G15:J15").ClearContents()
$ws2.Range("F15").Value = ""
$ws2.Range("A16:E16").ClearContents()
$ws2.Range("G16:J16").ClearContents()
$ws2.Range("F16").Value = ""
$ws2.Range("A17:E17").ClearContents()
$ws2.Range("G17:J17").ClearContents()
$ws2.Range("F17").Value = ""

# --- Selections / active sheet, matching the target sheetViews ---
$ws2.Activate()
$ws2.Range("H2:I14").Select()

$ws1.Activate()
$ws1.Range("A15:XFD114").Select()
